$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("001_LoginCorrectly")
$ws2 = $wb.Worksheets.Item("002_LoginIncorrectly")

# Rename Sheet3 and set it up with new data
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Name = "004_AccessoriesPage"

$ws3.Range("A1").Value = "TestCase"
$ws3.Range("A2").Value = "004-Test Accessories Page"

$ws3.Range("B1").Value = "Product Name"
$ws3.Range("B2").Value = "Magic Mouse"
$ws3.Range("B3").Value = "Apple TV"
$ws3.Range("B4").Value = "Sennheiser RS 120"
$ws3.Range("B5").Value = "Skullcandy PLYR 1 – Black"
$ws3.Range("B6").Value = "Apple 27 inch Thunderbolt Display"
$ws3.Range("B7").Value = "Asus MX239H 23-inch Widescreen AH"

$ws3.Range("C1").Value = "Product Prices"

$ws3.Range("C2:C7").NumberFormat = "@"
$ws3.Range("C2").Value = "`$150.00"
$ws3.Range("C3").Value = "`$80.00"
$ws3.Range("C4").Value = "`$50.00"
$ws3.Range("C5").Value = "`$110.00"
$ws3.Range("C6").Value = "`$764.00"
$ws3.Range("C7").Value = "`$199.00"

$ws3.Range("A1:C7").EntireColumn.AutoFit()

$ws3.PageSetup.PaperSize = 9
$ws3.PageSetup.Orientation = 1

$ws3.Range("C6").Select()
$ws3.Activate()

$ws1.Activate()
$ws1.Range("A14").Select()

$ws2.Activate()
$ws2.Range("A1:C3").Select()

$ws3.Activate()
